$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("E4").Value = "ht"
$ws2.Range("E5").Value = "ht"
$ws2.Range("E6").Value = "ht"
$ws2.Range("E7").Value = "ht"
$ws2.Range("H4").Value = "2016-09-03 00:34:32"
$ws2.Range("H5").Value = "2016-09-03 00:34:32"
$ws2.Range("H6").Value = "2016-09-03 00:34:32"
$ws2.Range("H7").Value = "2016-09-03 00:34:32"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("E4").Value = "ht"
$ws3.Range("E5").Value = "ht"
$ws3.Range("E6").Value = "ht"
$ws3.Range("E7").Value = "ht"
$ws3.Range("H4").Value = "2016-09-03 00:34:37"
$ws3.Range("H5").Value = "2016-09-03 00:34:37"
$ws3.Range("H6").Value = "2016-09-03 00:34:37"
$ws3.Range("H7").Value = "2016-09-03 00:34:37"

# "Latest HO Xliff Generate Date" on the Overview sheet mirrors the de-de
# "Latest Handoff Datetime" value for these rows, so it advances too.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("G4").Value = "2016-09-03 00:34:37"
$ws1.Range("G5").Value = "2016-09-03 00:34:37"
$ws1.Range("G6").Value = "2016-09-03 00:34:37"
$ws1.Range("G7").Value = "2016-09-03 00:34:37"
